# Updated symbol list (price/volume refresh + a BitrueCoin/MandalaExchangeToken row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, taken from the data refresh.
$updates = @(
    @{ Cell = "D2"; Value = "306.47" }
    @{ Cell = "E2"; Value = "-6.64%" }
    @{ Cell = "D3"; Value = "39.16" }
    @{ Cell = "E3"; Value = "-12.36%" }
    @{ Cell = "D4"; Value = "5.002" }
    @{ Cell = "E4"; Value = "-5.58%" }
    @{ Cell = "D5"; Value = "0.07726" }
    @{ Cell = "E5"; Value = "-7.70%" }
    @{ Cell = "D6"; Value = "4.279" }
    @{ Cell = "E6"; Value = "-3.38%" }
    @{ Cell = "D7"; Value = "1.594" }
    @{ Cell = "E7"; Value = "-18.20%" }
    @{ Cell = "D8"; Value = "0.9173" }
    @{ Cell = "E8"; Value = "-5.57%" }
    @{ Cell = "D9"; Value = "0.1008" }
    @{ Cell = "E9"; Value = "-9.03%" }
    @{ Cell = "D10"; Value = "0.1734" }
    @{ Cell = "E10"; Value = "-9.22%" }
    @{ Cell = "B11"; Value = "MandalaExchangeToken" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D11"; Value = "0.09027" }
    @{ Cell = "E11"; Value = "-6.89%" }
    @{ Cell = "B12"; Value = "BitrueCoin" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D12"; Value = "0.04484" }
    @{ Cell = "E12"; Value = "-2.54%" }
    @{ Cell = "D13"; Value = "7.046" }
    @{ Cell = "E13"; Value = "-15.71%" }
    @{ Cell = "D14"; Value = "0.1060" }
    @{ Cell = "E14"; Value = "-0.06%" }
    @{ Cell = "D15"; Value = "0.001256" }
    @{ Cell = "E15"; Value = "-2.65%" }
    @{ Cell = "E16"; Value = "-3.88%" }
    @{ Cell = "D17"; Value = "3.364" }
    @{ Cell = "E17"; Value = "-0.06%" }
    @{ Cell = "D19"; Value = "0.3366" }
    @{ Cell = "E19"; Value = "0.34%" }
    @{ Cell = "D20"; Value = "0.1363" }
    @{ Cell = "E20"; Value = "0.73%" }
    @{ Cell = "D21"; Value = "0.2864" }
    @{ Cell = "E21"; Value = "5.35%" }
    @{ Cell = "D22"; Value = "0.04147" }
    @{ Cell = "E22"; Value = "-0.69%" }
    @{ Cell = "D23"; Value = "0.001200" }
    @{ Cell = "E23"; Value = "-2.97%" }
    @{ Cell = "D24"; Value = "0.004080" }
    @{ Cell = "E24"; Value = "-8.27%" }
    @{ Cell = "E25"; Value = "-5.72%" }
    @{ Cell = "E26"; Value = "0.45%" }
    @{ Cell = "D38"; Value = "0.02340" }
    @{ Cell = "E38"; Value = "-13.87%" }
    @{ Cell = "D39"; Value = "0.05130" }
    @{ Cell = "E39"; Value = "-9.01%" }
    @{ Cell = "D40"; Value = "0.007945" }
    @{ Cell = "E40"; Value = "2.35%" }
    @{ Cell = "D41"; Value = "0.1328" }
    @{ Cell = "E41"; Value = "-5.91%" }
    @{ Cell = "D42"; Value = "0.007327" }
    @{ Cell = "E42"; Value = "0.10%" }
    @{ Cell = "D43"; Value = "0.001987" }
    @{ Cell = "E43"; Value = "-6.21%" }
    @{ Cell = "D44"; Value = "0.008019" }
    @{ Cell = "E44"; Value = "1.74%" }
    @{ Cell = "D45"; Value = "0.3310" }
    @{ Cell = "E45"; Value = "-5.62%" }
    @{ Cell = "E46"; Value = "-3.82%" }
    @{ Cell = "E47"; Value = "0.45%" }
    @{ Cell = "D48"; Value = "0.003397" }
    @{ Cell = "E48"; Value = "-2.61%" }
    @{ Cell = "E49"; Value = "16.67%" }
    @{ Cell = "E50"; Value = "0.45%" }
    @{ Cell = "E51"; Value = "0.45%" }
)

foreach ($u in $updates) {
    $col = ($u.Cell -replace '[0-9]+$', '')
    $rng = $ws.Range($u.Cell)
    if ($col -eq "D" -or $col -eq "E") {
        # Price / Volume(1h) columns are stored as plain text (e.g. "306.47", "-6.64%");
        # force text formatting so Excel does not silently coerce these into numbers/percentages,
        # then restore the default (unstyled) cell style so no stray formatting is introduced.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
